$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value  = 1.83
$ws.Range("I2").Value  = 5
$ws.Range("J2").Value  = 2.6
$ws.Range("L2").Value  = 5.5
$ws.Range("M2").Value  = 1.1
$ws.Range("N2").Value  = 7
$ws.Range("AE2").Value = 21
$ws.Range("AH2").Value = 23
$ws.Range("AN2").Value = 3.6
$ws.Range("AO2").Value = 10
$ws.Range("AX2").Value = 29
